$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B3:F6").BorderAround(1)
